$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.09665433333333333
$ws.Range("H2").Value = 0.289963
$ws.Range("I2").Value = 0.0006230336790718351
$ws.Range("J2").Value = 0.0006230336790718351
$ws.Range("M2").Value = 0.7435376666666667
$ws.Range("N2").Value = 2.230613
$ws.Range("O2").Value = 0.5505219265933909
$ws.Range("P2").Value = 0.5505219265933909
$ws.Range("Q2").Value = 0.07186613747988889
$ws.Range("R2").Value = 0.6467952373189999
$ws.Range("S2").Value = 0.0003429937013351951
$ws.Range("T2").Value = 0.0003429937013351951
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.09665433333333333
$ws.Range("H3").Value = 0.289963
$ws.Range("I3").Value = 0.0006230336790718351
$ws.Range("J3").Value = 0.0006230336790718351
$ws.Range("O3").Value = 0.08871012126664225
$ws.Range("P3").Value = 0.08871012126664224
$ws.Range("Q3").Value = 0.01158038120344444
$ws.Range("R3").Value = 0.104223430831
$ws.Range("S3").Value = 0.00005526939322366476
$ws.Range("T3").Value = 0.00005526939322366475
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.09665433333333333
$ws.Range("H4").Value = 0.289963
$ws.Range("I4").Value = 0.0006230336790718351
$ws.Range("J4").Value = 0.0006230336790718351
$ws.Range("M4").Value = 0.487255
$ws.Range("N4").Value = 1.461765
$ws.Range("O4").Value = 0.3607679521399669
$ws.Range("P4").Value = 0.3607679521399669
$ws.Range("Q4").Value = 0.04709530718833333
$ws.Range("R4").Value = 0.423857764695
$ws.Range("S4").Value = 0.0002247705845129753
$ws.Range("T4").Value = 0.0002247705845129753
$ws.Range("G5").Value = 154.8642143333334
$ws.Range("H5").Value = 464.5926430000001
$ws.Range("I5").Value = 0.9982544794956518
$ws.Range("J5").Value = 0.9982544794956519
$ws.Range("M5").Value = 0.7435376666666667
$ws.Range("N5").Value = 2.230613
$ws.Range("O5").Value = 0.5505219265933909
$ws.Range("P5").Value = 0.5505219265933909
$ws.Range("Q5").Value = 115.1473765755732
$ws.Range("R5").Value = 1036.326389180159
$ws.Range("S5").Value = 0.5495609792824289
$ws.Range("T5").Value = 0.549560979282429
$ws.Range("G6").Value = 154.8642143333334
$ws.Range("H6").Value = 464.5926430000001
$ws.Range("I6").Value = 0.9982544794956518
$ws.Range("J6").Value = 0.9982544794956519
$ws.Range("O6").Value = 0.08871012126664225
$ws.Range("P6").Value = 0.08871012126664224
$ws.Range("S6").Value = 0.08855527593102812
$ws.Range("T6").Value = 0.08855527593102812
$ws.Range("G7").Value = 154.8642143333334
$ws.Range("H7").Value = 464.5926430000001
$ws.Range("I7").Value = 0.9982544794956518
$ws.Range("J7").Value = 0.9982544794956519
$ws.Range("M7").Value = 0.487255
$ws.Range("N7").Value = 1.461765
$ws.Range("O7").Value = 0.3607679521399669
$ws.Range("P7").Value = 0.3607679521399669
$ws.Range("Q7").Value = 75.45836275498834
$ws.Range("R7").Value = 679.1252647948951
$ws.Range("S7").Value = 0.3601382242821949
$ws.Range("T7").Value = 0.3601382242821949
$ws.Range("G8").Value = 0.174137
$ws.Range("H8").Value = 0.522411
$ws.Range("I8").Value = 0.001122486825276316
$ws.Range("J8").Value = 0.001122486825276316
$ws.Range("M8").Value = 0.7435376666666667
$ws.Range("N8").Value = 2.230613
$ws.Range("O8").Value = 0.5505219265933909
$ws.Range("P8").Value = 0.5505219265933909
$ws.Range("Q8").Value = 0.1294774186603333
$ws.Range("R8").Value = 1.165296767943
$ws.Range("S8").Value = 0.0006179536096268163
$ws.Range("T8").Value = 0.0006179536096268166
$ws.Range("G9").Value = 0.174137
$ws.Range("H9").Value = 0.522411
$ws.Range("I9").Value = 0.001122486825276316
$ws.Range("J9").Value = 0.001122486825276316
$ws.Range("O9").Value = 0.08871012126664225
$ws.Range("P9").Value = 0.08871012126664224
$ws.Range("Q9").Value = 0.02086376028966667
$ws.Range("R9").Value = 0.187773842607
$ws.Range("S9").Value = 0.00009957594239047025
$ws.Range("T9").Value = 0.00009957594239047026
$ws.Range("G10").Value = 0.174137
$ws.Range("H10").Value = 0.522411
$ws.Range("I10").Value = 0.001122486825276316
$ws.Range("J10").Value = 0.001122486825276316
$ws.Range("M10").Value = 0.487255
$ws.Range("N10").Value = 1.461765
$ws.Range("O10").Value = 0.3607679521399669
$ws.Range("P10").Value = 0.3607679521399669
$ws.Range("Q10").Value = 0.08484912393499999
$ws.Range("R10").Value = 0.7636421154149999
$ws.Range("S10").Value = 0.0004049572732590293
$ws.Range("T10").Value = 0.0004049572732590294
